$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Header row (column order/content unchanged, but make sure B1 is present) ---
$ws.Range("A1").Value = "Senario"
$ws.Range("B1").Value = "user_skill_id"
$ws.Range("C1").Value = "user_id"
$ws.Range("D1").Value = "Skill_id"
$ws.Range("E1").Value = "months_of_exp"
$ws.Range("F1").Value = "StatusCode"
$ws.Range("G1").Value = "StatusMessage"

# --- Clear the leftover "Docs-Calibri" formatting carried on A2:A3 from the old rows ---
$ws.Range("A2:A8").ClearFormats()

# --- Row 2 ---
$ws.Range("A2").Value = "To map new user and skill"
$ws.Range("B2").Value = "US15"
$ws.Range("C2").Value = "U04"
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 201
$ws.Range("G2").Value = "Successfully Created"

# --- Row 3 ---
$ws.Range("A3").Value = "To map to existing user And Skills"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = "U04"
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 400
$ws.Range("G3").Value = "Failed to create as UserSkillMap already exists"

# --- Row 4 ---
$ws.Range("A4").Value = "To map new user and skill , with skill id as alpha numeric"
$ws.Range("C4").Value = "U07"
$ws.Range("D4").Value = "A12"
$ws.Range("E4").Value = 12
$ws.Range("F4").Value = 400
$ws.Range("G4").Value = "Failed to create due to invalid data"

# --- Row 5 ---
$ws.Range("A5").Value = "To map new user and skill , with skill id as null"
$ws.Range("C5").Value = "U07"
$ws.Range("E5").Value = 12
$ws.Range("F5").Value = 400
$ws.Range("G5").Value = "Failed to create due to invalid skill Id"

# --- Row 6 ---
$ws.Range("A6").Value = "To map new user and skill , with user id as null"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 11
$ws.Range("F6").Value = 400
$ws.Range("G6").Value = "Failed to create due to invalid data"

# --- Row 7 ---
$ws.Range("A7").Value = "To map new user and skill , with months of experience as alpha numeric"
$ws.Range("C7").Value = "U07"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "a11"
$ws.Range("F7").Value = 400
$ws.Range("G7").Value = "Failed to create due to invalid data"

# --- Row 8 ---
$ws.Range("A8").Value = "To map new user and skill , with months of experience as null"
$ws.Range("C8").Value = "U07"
$ws.Range("D8").Value = 3
$ws.Range("F8").Value = 400
$ws.Range("G8").Value = "Failed to create due to invalid data"

# --- Touch row 9 (stays blank) so the sheet's used range extends to it ---
$ws.Range("A9:G9").Interior.Pattern = -4142

# --- Column widths (best effort match to the target character widths) ---
$ws.Columns.Item(1).ColumnWidth = 46.25
$ws.Columns.Item(2).ColumnWidth = 19.42

Write-Host "edit complete"
